# supergroup.xlsx ("EERV Région 2 / 2030") - first pass at fixing the
# parish-group ids (IdxPar, column B): the sheet had been cloned from the
# "4020" group and still carried its id (4020000000) on every row instead
# of the correct "2030" (2030000000) id for this region.
#
# While at it: drop the one clearly-bogus row (row 35) whose IdxSG id
# (100000000558) has too many digits to be real and whose label
# ("ACTIVITES NON CLASSEES") goes away with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the parish id in column B for every data row (2..34): 4020000000 -> 2030000000
$ws.Range("B2:B34").Value = 2030000000

# Row 35 is bogus (impossible 12-digit id) - wipe it out, keeping the (now empty) row.
$ws.Range("B35:D35").ClearContents()

# Leave the selection where the author ended up after editing.
$ws.Range("C31").Select() | Out-Null
